$wb = $excel.ActiveWorkbook
$wsTags = $wb.Worksheets.Item("tags")
$wsPackages = $wb.Worksheets.Item("packages")

# --- Add the new "sam_biosharing" tag as row 4 of the "tags" sheet ---
# Copy the formatting of the last existing row (row 3) down into row 4 first,
# so the new row inherits the same per-cell styling (F column keeps its
# existing "system" style) without Excel minting brand-new style records.
$wsTags.Range("A3:F3").Copy()
$wsTags.Range("A4:F4").PasteSpecial(-4122)

$wsTags.Range("A4").Value = "sam_biosharing"
$wsTags.Range("D4").Value = "Reference"
$wsTags.Range("B4").Value = "http://www.biosharing.org/bsg-000261"
$wsTags.Range("C4").Value = "http://www.biosharing.org/bsg-000261"
$wsTags.Range("E4").Value = "system"
$wsTags.Range("F4").Value = "http://molgenis.org/biobankconnect/link"

# The wider new row now drives the visible columns, so fit them to content
# (mirrors the column widths Excel would compute/persist on save).
$wsTags.Range("A1:F4").EntireColumn.AutoFit()

# --- Reference the new tag from the SAM package's "tags" column ---
$wsPackages.Range("D2").Value = "sam_doc,sam_pub1,sam_biosharing"

# --- Update selection / active sheet to match the saved view state ---
$wsTags.Activate()
$wsTags.Range("A4").Select()

$wsPackages.Activate()
$wsPackages.Range("D2").Select()
